$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/25/2025  Through  8/31/2025"

# --- Cells that change between numeric and text-marker type (need format fix-up) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$ws.Range("I31").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2

$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 150

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# --- Plain numeric value updates ---
# Row 14
$ws.Range("M14").Value = -84.210526315789

# Row 15
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 6
$ws.Range("I15").Value = 21
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = -16
$ws.Range("L15").Value = -27.586206896551
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = -66.129032258064

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 167
$ws.Range("K16").Value = 19.760479041916
$ws.Range("L16").Value = 12.359550561797
$ws.Range("M16").Value = -32.885906040268
$ws.Range("N16").Value = -86.675549633577

# Row 17
$ws.Range("C17").Value = 11
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 49
$ws.Range("H17").Value = 12.244897959183
$ws.Range("I17").Value = 564
$ws.Range("J17").Value = 463
$ws.Range("K17").Value = 21.814254859611
$ws.Range("L17").Value = 15.811088295687
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -29.145728643216

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 166.666666666667
$ws.Range("I18").Value = 130
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = 36.842105263157
$ws.Range("L18").Value = 54.761904761904
$ws.Range("M18").Value = -33.673469387755
$ws.Range("N18").Value = -75.190839694656

# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 56.521739130434
$ws.Range("I19").Value = 313
$ws.Range("J19").Value = 248
$ws.Range("K19").Value = 26.209677419354
$ws.Range("L19").Value = 36.681222707423
$ws.Range("M19").Value = 7.931034482758
$ws.Range("N19").Value = -49.188311688311

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -26.315789473684
$ws.Range("I20").Value = 103
$ws.Range("J20").Value = 106
$ws.Range("K20").Value = -2.830188679245
$ws.Range("L20").Value = -4.629629629629
$ws.Range("M20").Value = -1.904761904761
$ws.Range("N20").Value = -75.180722891566

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = 22.950819672131
$ws.Range("I21").Value = 1334
$ws.Range("J21").Value = 1118
$ws.Range("K21").Value = 19.320214669051
$ws.Range("L21").Value = 19.000892060660
$ws.Range("M21").Value = -1.258327165062
$ws.Range("N21").Value = -66.321635950517

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = -33.333333333333
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -62.5
$ws.Range("I22").Value = 35
$ws.Range("J22").Value = 47
$ws.Range("K22").Value = -25.531914893617
$ws.Range("L22").Value = 12.903225806451
$ws.Range("M22").Value = 6.060606060606

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 200
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 219
$ws.Range("J23").Value = 245
$ws.Range("K23").Value = -10.612244897959
$ws.Range("L23").Value = -10.612244897959
$ws.Range("M23").Value = 36.875

# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 13.333333333333
$ws.Range("I24").Value = 854
$ws.Range("J24").Value = 767
$ws.Range("K24").Value = 11.342894393741
$ws.Range("L24").Value = 12.664907651715
$ws.Range("M24").Value = 29.590288315629

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 46.666666666666
$ws.Range("I25").Value = 226
$ws.Range("J25").Value = 196
$ws.Range("K25").Value = 15.306122448979
$ws.Range("L25").Value = 10.784313725490

# Row 26
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 6.25
$ws.Range("F26").Value = 88
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 12.820512820512
$ws.Range("I26").Value = 680
$ws.Range("J26").Value = 572
$ws.Range("K26").Value = 18.881118881118
$ws.Range("L26").Value = 3.186646433990
$ws.Range("M26").Value = -24.611973392461

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("I27").Value = 27
$ws.Range("J27").Value = 34
$ws.Range("K27").Value = -20.588235294117
$ws.Range("L27").Value = -20.588235294117

# Row 28
$ws.Range("C28").Value = 5
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 325
$ws.Range("I28").Value = 77
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = 71.111111111111
$ws.Range("L28").Value = 28.333333333333

# Row 29
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("M29").Value = -66.25
$ws.Range("N29").Value = -87.892376681614

# Row 30
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("M30").Value = -66.153846153846
$ws.Range("N30").Value = -89.320388349514
